$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card16")

# --- Row 19: fill in the previously-blank B:K cells with the literal text "nan" ---
$ws.Range("B19:K19").Value = "nan"

# --- Row 20: brand new service-log entry ---
# Column A holds "16" as text (matches the rest of column A on this sheet).
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "16"
$ws.Range("A20").Style = "Normal"

# Columns B:K are left blank for this entry (mirrors the other rows' pattern).
$ws.Range("B20:K20").NumberFormat = "@"
$ws.Range("B20:K20").Style = "Normal"

$ws.Range("L20").Value = "21\10\2025"
$ws.Range("M20").Value = "8670 h   696 t"
$ws.Range("N20").Value = "تم تغيير الجرائد الخلفيه (1_5_8) ومعيار المكنه"
$ws.Range("O20").Value = "الخبير"
